# add 2d act camera
# Update the camera offset position/rotation for the first scene row
# (villageScene), and leave the active selection on the edited cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = "0,4.2,5.5"
$ws.Range("K2").Value = "25,180"

$ws.Range("K2").Select()
